# Daily attendance processing - 2025-12-04 17:54:18
# Reorders the "Recorded By" (column G) entries: for any cell whose value is a
# comma-separated list of recorders, the first two entries are swapped
# (e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ","
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }
        if ($parts.Length -ge 2) {
            $tmp = $parts[0]
            $parts[0] = $parts[1]
            $parts[1] = $tmp
            $cell.Value2 = [string]::Join(", ", $parts)
        }
    }
}
